$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 248.1910118546626
